# Generate Report for Handoff
#
# The localization status report is being refreshed: the "handed back"
# status / timestamps from the previous sync are replaced with a fresh
# "ready for handoff" status + new generation timestamps, and the
# (now shorter) status column is narrowed to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refresh the handoff generation timestamps ---
$overview.Range("G2").Value = "2016-09-02 21:08:44"
$zhcn.Range("H2").Value     = "2016-09-02 21:08:39"
$dede.Range("H2").Value     = "2016-09-02 21:08:44"

# --- Narrow the status columns now that the text is shorter ---
$overview.Columns.Item(5).ColumnWidth = 16.25
$overview.Columns.Item(6).ColumnWidth = 16.25
$zhcn.Columns.Item(3).ColumnWidth     = 16.25
$dede.Columns.Item(3).ColumnWidth     = 16.25
